$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 401 (weekly update pushes existing rows down)
$ws.Rows.Item(401).Insert()
$ws.Rows.Item(401).Insert()
$ws.Rows.Item(401).Insert()

# Populate the 3 new rows with the new week of data (Provincia de Quillota, 2022-02-18)
$ws.Range("A401").Value = 3
$ws.Range("B401").Value = "Femacal de La Calera"
$ws.Range("C401").Value = "Coquimbo"
$ws.Range("D401").Value = 44610
$ws.Range("E401").Value = 5
$ws.Range("F401").Value = 100112027
$ws.Range("G401").Value = "Melón"
$ws.Range("H401").Value = "Tuna"
$ws.Range("I401").Value = "Extra"
$ws.Range("J401").Value = 120
$ws.Range("K401").Value = 1100
$ws.Range("L401").Value = 1100
$ws.Range("M401").Value = 1100
$ws.Range("N401").Value = "`$/unidad"
$ws.Range("O401").Value = "Provincia de Quillota"
$ws.Range("P401").Value = 1100
$ws.Range("Q401").Value = 1
$ws.Range("R401").Value = "Hortaliza"
$ws.Range("A402").Value = 3
$ws.Range("B402").Value = "Femacal de La Calera"
$ws.Range("C402").Value = "Coquimbo"
$ws.Range("D402").Value = 44610
$ws.Range("E402").Value = 5
$ws.Range("F402").Value = 100112027
$ws.Range("G402").Value = "Melón"
$ws.Range("H402").Value = "Tuna"
$ws.Range("I402").Value = "Primera"
$ws.Range("J402").Value = 210
$ws.Range("K402").Value = 700
$ws.Range("L402").Value = 700
$ws.Range("M402").Value = 700
$ws.Range("N402").Value = "`$/unidad"
$ws.Range("O402").Value = "Provincia de Quillota"
$ws.Range("P402").Value = 700
$ws.Range("Q402").Value = 1
$ws.Range("R402").Value = "Hortaliza"
$ws.Range("A403").Value = 3
$ws.Range("B403").Value = "Femacal de La Calera"
$ws.Range("C403").Value = "Coquimbo"
$ws.Range("D403").Value = 44610
$ws.Range("E403").Value = 5
$ws.Range("F403").Value = 100112027
$ws.Range("G403").Value = "Melón"
$ws.Range("H403").Value = "Tuna"
$ws.Range("I403").Value = "Segunda"
$ws.Range("J403").Value = 200
$ws.Range("K403").Value = 500
$ws.Range("L403").Value = 500
$ws.Range("M403").Value = 500
$ws.Range("N403").Value = "`$/unidad"
$ws.Range("O403").Value = "Provincia de Quillota"
$ws.Range("P403").Value = 500
$ws.Range("Q403").Value = 1
$ws.Range("R403").Value = "Hortaliza"
